# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
# Reorder the "Periodo Mora" column (E16:E22) so the newest periods are listed
# first (2308 down to 2302) and keep the "Valor Mora" date/amount cell (F)
# attached to the correct row after the reorder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order of periods for rows 16-22 (previously 2302..2308 top to bottom)
$ws.Range("E16").Value = "2308"
$ws.Range("E17").Value = "2307"
$ws.Range("E18").Value = "2306"
$ws.Range("E19").Value = "2305"
$ws.Range("E20").Value = "2304"
$ws.Range("E21").Value = "2303"
$ws.Range("E22").Value = "2302"

# The value that used to sit on the "2308" row (46400) now belongs to row 16,
# and the value that used to sit on the "2302" row (44854) now belongs to row 22.
$ws.Range("F16").Value = 46400
$ws.Range("F22").Value = 44854
